# Apply updated symbol list values (Price and Volume(1h) columns)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "327.52"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.29%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "44.42"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.13%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.121"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-7.07%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08374"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "3.85%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.443"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.31%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.937"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-5.55%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9744"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.41%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.521"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-3.82%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1124"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.45%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1902"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.22%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09672"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-2.42%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.04623"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-2.04%"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.36%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001290"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "2.43%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005802"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-5.39%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.404"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.96%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.3361"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.45%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.887"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-13.08%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1363"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-2.27%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.01%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04166"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.59%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001239"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-5.55%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004424"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "2.05%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001306"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.62%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0002985"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-20.45%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02722"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "3.44%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05632"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "0.55%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007849"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.78%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1413"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.74%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007367"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.69%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002118"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "6.29%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007926"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-9.23%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3500"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006911"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-3.01%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000753"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.07%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003508"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-1.80%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.003538"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "39.89%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002109"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.07%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002008"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.07%"
